$wb = $excel.ActiveWorkbook

$notas = $wb.Worksheets.Item("Notas")
$tarefas = $wb.Worksheets.Item("Tarefas")

# Fix Lucas' score for Sprint 2 (was 0, now 4)
$notas.Range("E4").Value = 4

# Add note about Lucas delaying the sprint delivery
$notas.Range("H4").Value = "Lucas atrasou a entrega da Sprint"

# Add the Sprint 3 row of data
$notas.Range("A5").Value = 3
$notas.Range("B5").Value = 8
$notas.Range("C5").Value = 10
$notas.Range("D5").Value = 8
$notas.Range("E5").Value = 8
$notas.Range("F5").Value = 10

# Add note about Marina/Camilli workload
$notas.Range("H5").Value = "Marina teve mais trabalho que os demais, Camilli ajudou Marina"

# Notas' own cursor moves to A6 (just past the new Sprint-3 row) before
# the workbook focus moves away from it
$notas.Range("A6").Select()

# Make Tarefas the active sheet/tab
$tarefas.Activate()

$wb.Save()
